$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add PR reference hyperlink for Athmika Bhat (row 12) ---
$prUrl  = "https://github.com/dhavalkeerthi/MRIInterns2026A/pull/24"
$prText = "Add Athmika.Bhat.txt file by AthmikaU · Pull Request #24 · dhavalkeerthi/MRIInterns2026A"

$target = $ws.Range("B12")

# Create the hyperlink relationship (this also stamps a "Hyperlink" style on
# the cell as a side effect of Excel's own Insert Hyperlink behaviour).
$ws.Hyperlinks.Add($target, $prUrl, "", "", $prUrl)

# The visible cell text is the descriptive PR title, not the bare URL.
$target.Value = $prText

# Make sure the cell uses the same "Hyperlink" look as the existing B4 link.
$target.Style = "Hyperlink"

# --- Widen column B to fit the new, longer text and drop the old bestFit flag ---
$ws.Columns("B").ColumnWidth = 73.8

# --- Update the active selection left behind by the edit ---
$ws.Range("C12").Select() | Out-Null
